$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 630
$ws.Range("F3").Value = 288
$ws.Range("F6").Value = 447
$ws.Range("G6").Value = 36
$ws.Range("F8").Value = 197
$ws.Range("F10").Value = 272
$ws.Range("F11").Value = 7111
$ws.Range("F13").Value = 59
$ws.Range("F14").Value = 274
$ws.Range("F16").Value = 566
$ws.Range("F17").Value = 391
$ws.Range("F18").Value = 431
$ws.Range("F22").Value = 6
$ws.Range("F23").Value = 30
$ws.Range("F24").Value = 95
$ws.Range("F27").Value = 110
$ws.Range("F28").Value = 345
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 1052
$ws.Range("F32").Value = 65
$ws.Range("F33").Value = 2034
$ws.Range("F34").Value = 560
$ws.Range("F35").Value = 2
$ws.Range("F37").Value = 34
$ws.Range("F38").Value = 547
$ws.Range("F39").Value = 10

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 279
$ws.Range("F4").Value = 58
$ws.Range("F5").Value = 284
$ws.Range("F8").Value = 49

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 347

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 347
$ws.Range("F3").Value = 630
$ws.Range("F4").Value = 288
$ws.Range("F8").Value = 447
$ws.Range("G8").Value = 36
$ws.Range("F10").Value = 197
$ws.Range("F12").Value = 272
$ws.Range("F13").Value = 7111
$ws.Range("F15").Value = 59
$ws.Range("F16").Value = 279
$ws.Range("F17").Value = 274
$ws.Range("F19").Value = 566
$ws.Range("F20").Value = 391
$ws.Range("F21").Value = 431
$ws.Range("F22").Value = 58
$ws.Range("F25").Value = 284
$ws.Range("F29").Value = 6
$ws.Range("F30").Value = 30
$ws.Range("F31").Value = 95
$ws.Range("F32").Value = 49
$ws.Range("F37").Value = 110
$ws.Range("F38").Value = 345
$ws.Range("F39").Value = 4
$ws.Range("F40").Value = 1052
$ws.Range("F42").Value = 65
$ws.Range("F43").Value = 2034
$ws.Range("F44").Value = 560
$ws.Range("F45").Value = 2
$ws.Range("F47").Value = 34
$ws.Range("F48").Value = 547
$ws.Range("F49").Value = 10
